# "changed correlation factor pump"
# Updates the pump correlation-factor values (column P, feed_rate_mL_to_g)
# from 1.14 to 1.164 for the active runs, adjusts the related unit label
# and the row-7 "mL/%/min" calibration value, and moves the active
# selection to reflect where the editor was working when the file was
# saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K2: unit header for mX0 changed from "g/L" to "g" ---
$ws.Range("K2").Value = "g"

# --- P3:P6: correlation factor pump updated 1.14 -> 1.164 ---
$ws.Range("P3:P6").Value = 1.164

# --- Row 7: O7 calibration value changed, P7 correlation factor updated
#     and its bottom border removed (it is no longer the last bordered row) ---
$ws.Range("O7").Value = 0.057
$ws.Range("P7").Value = 1.164
$ws.Range("P7").Borders.Item(9).LineStyle = -4142

# --- Selection moved to where the editor was last working ---
$ws.Range("O8").Select()
